# Update "想去人数" (want-to-go count) figures on the "展览" sheet
# and the combined "全部类型" sheet, as published to gh-pages at 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 125
$wsExpo.Range("F4").Value = 163
$wsExpo.Range("F5").Value = 3182
$wsExpo.Range("F6").Value = 318
$wsExpo.Range("F7").Value = 10

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 125
$wsAll.Range("F4").Value = 163
$wsAll.Range("F5").Value = 3182
$wsAll.Range("F6").Value = 318
$wsAll.Range("F9").Value = 10
